$wb = $excel.ActiveWorkbook

# "SoCDTtiNTY-psgr" sheet: update the passenger share-that-is-new formula
# from 1/20 to 1/17 (newest US data), keeping the existing shared-formula
# layout (B2 stand-alone, C2:H2 shared group).
$ws = $wb.Worksheets.Item("SoCDTtiNTY-psgr")
$ws.Range("B2").Formula = "=1/17"
$ws.Range("C2:H2").Formula = "=1/17"

# Move the active selection to F9, matching where the editor left off.
$ws.Range("F9").Select() | Out-Null
